$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells that hold numeric-looking strings so Excel
# does not silently convert them to numbers (matches how these price/volume
# cells were already stored as inline text strings in the workbook).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.386.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.721.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4871"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06183"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.16"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06958"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.43"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.538"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5961"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.396.10"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007204"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.940.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.437"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.473"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.101"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.722"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.925"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08000"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.665"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04488"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.604"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9975"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6233"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9411"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.387"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.946"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01471"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.89"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.321"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3829"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.825"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1162"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.723"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.226"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.74"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.41%  "
